$wb = $excel.ActiveWorkbook

# --- Sheet "financal and other parameteres" (sheet index 3): drop the
# demand_th / heat_storage / radiation / temp columns -------------------
$ws3 = $wb.Worksheets.Item("financal and other parameteres")

$ws3.Range("F1").EntireColumn.Delete()
$ws3.Range("E1").EntireColumn.Delete()
$ws3.Range("C1").EntireColumn.Delete()
$ws3.Range("B1").EntireColumn.Delete()

$ws3.Range("C12").Select()

# --- Active tab moves from "Heat Storage" back to "Heat Generators" ----
$ws1 = $wb.Worksheets.Item("Heat Generators")
$ws1.Activate()
$ws1.Range("C12").Select()

$wb.Save()
